# Kostendifferenzierung eingebaut, bessere Werte in Knotentabelle eingefügt

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Bemessungsleistung PK" (D) values for rows 2..10
$newD = @{
    2  = 1000000
    3  = 2000000
    4  = 700000
    5  = 300000
    6  = 2000000
    7  = 100000
    8  = 1500000
    9  = 750000
    10 = 400000
}

foreach ($row in 2..10) {
    $ws.Cells.Item($row, 4).Value = $newD[$row]
    $dAddr = "D" + $row
    $ws.Cells.Item($row, 5).Formula = "=(" + $dAddr + "*12)/60+(" + $dAddr + "/10000)*500000/20"
}

# Sheet view changes
$window = $excel.ActiveWindow
$window.Zoom = 185
$window.ScrollRow = 1
$window.ScrollColumn = 1
$ws.Range("D10").Select()
